# Adding two new Test Cases for Notifications (F10, F11) to the "Test Cases" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11: copy formatting from the last existing data row (row 10), then overwrite values.
$ws.Range("A10:E10").Copy($ws.Range("A11:E11"))
$ws.Range("A11").Value = "TestCase_F10"
$ws.Range("B11").Value = "OPQA-217"
$ws.Range("C11").Value = "Verify that user receives a notification when someone comments on an post contained in his watchlist"
$ws.Range("D11").Value = "Y"
$ws.Range("E11").Value = "PASS"

# Row 12: same pattern for the second new test case.
$ws.Range("A10:E10").Copy($ws.Range("A12:E12"))
$ws.Range("A12").Value = "TestCase_F11"
$ws.Range("B12").Value = "OPQA-218"
$ws.Range("C12").Value = "Verify that user receives a notification if someone likes his comment on a post"
$ws.Range("D12").Value = "Y"
$ws.Range("E12").Value = "PASS"

# Match the selection recorded in the saved workbook.
$ws.Range("D9").Select() | Out-Null
